$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add I0 and IF headers, copying the formatting of the
# existing header cell H1 (style index 1: bold, bordered, centered) ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-24: populate column I (I0) and column J (IF) ---
$jValues = @{
    2  = 5
    3  = 6
    4  = 5
    5  = 7
    6  = 8
    7  = 6
    8  = 7
    9  = 6
    10 = 7
    11 = 5
    12 = 5
    13 = 7
    14 = 6
    15 = 2
    16 = 5
    17 = 5
    18 = 6
    19 = 5
    20 = 4
    21 = 4
    22 = 6
    23 = 3
    24 = 2
}

$iValues = @{
    22 = 4
}

for ($row = 2; $row -le 24; $row++) {
    $iVal = 1
    if ($iValues.ContainsKey($row)) {
        $iVal = $iValues[$row]
    }
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jValues[$row]
}
